$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.121.80"
$ws.Range("E2").Value = "  +0.39%  "

# Row 3
$ws.Range("D3").Value = "3.116.96"
$ws.Range("E3").Value = "  +0.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'579.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "

# Row 6
$ws.Range("D6").Value = "'174.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.24%  "

# Row 9
$ws.Range("D9").Value = "'6.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("E10").Value = "  +0.04%  "

# Row 11
$ws.Range("D11").Value = "'0.479"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12
$ws.Range("D12").Value = "'0.0000248"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("D13").Value = "'37.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.57%  "

# Row 14
$ws.Range("E14").Value = "  -1.59%  "

# Row 15
$ws.Range("D15").Value = "3.635.48"

# Row 16
$ws.Range("D16").Value = "67.108.43"
$ws.Range("E16").Value = "  +0.42%  "

# Row 17
$ws.Range("D17").Value = "'7.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "

# Row 18
$ws.Range("D18").Value = "3.119.77"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19
$ws.Range("D19").Value = "'16.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.62%  "

# Row 20
$ws.Range("D20").Value = "'492.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.75%  "

# Row 21
$ws.Range("B21").Value = "Polygon"
$ws.Range("C21").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D21").Value = "'0.706"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'7.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.97%  "

# Row 23
$ws.Range("D23").Value = "'84.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "

# Row 24
$ws.Range("D24").Value = "'13.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").Value = "'2.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.44%  "

# Row 26
$ws.Range("D26").Value = "'10.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.29%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("E28").Value = "  -0.69%  "

# Row 29
$ws.Range("D29").Value = "'2.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("E30").Value = "  -0.36%  "

# Row 31
$ws.Range("D31").Value = "'28.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.84%  "

# Row 32
$ws.Range("E32").Value = "  -0.73%  "

# Row 33
$ws.Range("D33").Value = "0.0₃0947"
$ws.Range("E33").Value = "  -5.31%  "

# Row 34
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").Value = "'5.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "

# Row 36
$ws.Range("E36").Value = "  -1.39%  "

# Row 37
$ws.Range("D37").Value = "'47.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.69%  "

# Row 38
$ws.Range("E38").Value = "  -2.58%  "

# Row 39
$ws.Range("D39").Value = "'0.311"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.52%  "

# Row 40
$ws.Range("D40").Value = "'0.123"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.41%  "

# Row 41
$ws.Range("D41").Value = "'8.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.74%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'387.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.821.96"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44
$ws.Range("E44").Value = "  -6.62%  "

# Row 45
$ws.Range("D45").Value = "'0.0352"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.29%  "

# Row 46
$ws.Range("D46").Value = "'135.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.24%  "

# Row 47
$ws.Range("E47").Value = "  +0.00%  "

# Row 48
$ws.Range("E48").Value = "  +0.15%  "

# Row 49
$ws.Range("E49").Value = "  -0.98%  "

# Row 50
$ws.Range("E50").Value = "  -0.73%  "

# Row 51
$ws.Range("D51").Value = "'6.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.57%  "
